$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for A2:H31 (timestamp, label, ax, ay, az, gx, gy, gz)
$rows = @(
  @(0, "walkingToRunning", -1.350502490997314,-5.659398555755615,-3.780858278274536,-0.2807844198983286,-0.8123789481047923,-0.09823161813205905),
  @(100, "walkingToRunning", 6.421082973480225,-10.34209537506104,-8.60602855682373,-0.1417398627462046,-0.1100757933896112,2.834019057966533),
  @(200, "walkingToRunning", -8.985816955566406,-4.653406143188477,-4.82445764541626,-1.13705500035451,2.739093329500554,-2.630493768330314),
  @(300, "walkingToRunning", 8.638599395751953,-36.5815658569336,4.202666282653809,2.9049899660308,-3.484631538391146,-2.321273335095107),
  @(400, "walkingToRunning", 0.9458236694335938,-4.780321598052979,5.630753517150879,-3.415270102435171,0.5314338042818449,4.029043727907666),
  @(500, "walkingToRunning", -5.881309509277344,-30.26239395141602,12.90286827087402,-9.314493606830482,-0.2480275487077682,3.206530965607738),
  @(600, "walkingToRunning", -47.72857666015625,-15.63389015197754,-17.3016414642334,2.801971287562949,-2.823644173042532,0.2525709908584024),
  @(700, "walkingToRunning", -12.29190921783447,-4.77882719039917,-10.05650520324707,4.358451415752542,-4.83039114393037,0.2956774553348279),
  @(800, "walkingToRunning", 21.71150207519531,-27.2324275970459,-6.215962886810303,7.929245866578198,4.034436612293629,0.2522587591204151),
  @(900, "walkingToRunning", 32.41067504882812,-4.009072303771973,21.43018341064453,-0.834473840121543,4.45597813047212,2.253706974202184),
  @(1000, "walkingToRunning", -4.226425170898438,-4.756555557250977,-0.5562095642089844,-5.379544188236389,2.00216909318133,4.38414403076831),
  @(1100, "walkingToRunning", 2.544495105743408,-37.85799407958984,14.00184154510498,-1.495801197043751,11.70056653022773,0.6984025085794245),
  @(1200, "walkingToRunning", 15.75080490112305,-12.57144546508789,-16.15233612060547,2.489601142961419,1.782228860361755,3.894388745057194),
  @(1300, "walkingToRunning", 77.36283111572266,-23.39109230041504,3.477614402770996,7.727643826912207,6.834787146798529,1.897908695812871),
  @(1400, "walkingToRunning", -13.44992446899414,-11.07156848907471,10.24287509918213,0.1587105783922702,5.019907260763223,1.782309127265012),
  @(1500, "walkingToRunning", -2.557040929794312,-6.614107131958008,-1.701248645782471,-5.217137336730959,3.712993112103693,5.708451303942454),
  @(1600, "walkingToRunning", -82.30638885498047,-41.48257446289063,-24.42162132263184,1.931081188136128,1.861609384931368,-5.062639532418077),
  @(1700, "walkingToRunning", 12.51860809326172,-2.532787799835205,-16.1464672088623,3.396664668773761,-5.718913538702679,-1.399837789864275),
  @(1800, "walkingToRunning", 25.23304176330566,-12.49789237976074,10.0257043838501,6.86892145666587,-1.869332508794124,1.079519285724066),
  @(1900, "walkingToRunning", 25.69987106323243,-2.908939123153687,11.61045837402344,1.767213952952357,1.710883551630483,1.963973425585661),
  @(2000, "walkingToRunning", 7.486310958862305,-1.500723719596863,4.869987487792969,-3.255076243959633,0.04516811000888521,5.548670966049739),
  @(2100, "walkingToRunning", -9.412309646606444,-35.61527252197266,2.476359844207764,-0.7715866319064268,8.701057298430134,-0.006967774752950184),
  @(2200, "walkingToRunning", 10.05931091308594,-4.656819343566895,-13.41165542602539,1.181180877932187,-0.4844702030050883,1.036063860202677),
  @(2300, "walkingToRunning", -0.115725427865982,-15.31198501586914,-12.16327095031738,5.493912943478295,-2.113819895119489,1.030801805956604),
  @(2400, "walkingToRunning", -12.95152854919434,-7.215466022491455,2.712479591369629,2.463171288884921,2.284831554725252,0.6873917538544219),
  @(2500, "walkingToRunning", -1.57819402217865,4.682579040527344,-4.256365299224854,-3.070133529860395,3.359637453638269,3.026170623713525),
  @(2600, "walkingToRunning", -8.448655128479004,-33.68074035644531,-12.51845359802246,-2.565556277488823,10.97180558895251,1.010093347779601),
  @(2700, "walkingToRunning", 71.44110107421875,6.458947658538818,-16.25710296630859,-0.5303593524571014,2.632135161037951,0.3993937804781349),
  @(2800, "walkingToRunning", 7.958407878875732,-4.263494491577148,-17.78672790527344,4.747484356164973,-6.67747688293454,-1.334976077079803),
  @(2900, "walkingToRunning", -30.355712890625,-23.57818984985352,8.899335861206055,2.688300560260526,-1.088815121815081,-1.408782399933858)
)

$data = New-Object 'object[,]' $rows.Count,8
for ($i = 0; $i -lt $rows.Count; $i++) {
  for ($j = 0; $j -lt 8; $j++) {
    $data[$i, $j] = $rows[$i][$j]
  }
}

$ws.Range("A2:H31").Value = $data
